$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delta Smelt (row 8, Adult) thresholds
$ws.Range("D8").Value = 19
$ws.Range("E8").Value = 22

# Overbite Clam (row 53, Undifferentiated) thresholds
$ws.Range("D53").Value = 28
$ws.Range("E53").Value = 28
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "NA"

$ws.Range("D53").Select()
